$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.260.87'
$ws.Range("E2").Value = '  -5.00%  '
$ws.Range("D3").Value = '2.237.66'
$ws.Range("E3").Value = '  -5.93%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.26'
$ws.Range("E5").Value = '  +0.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.66'
$ws.Range("E6").Value = '  -8.62%  '
$ws.Range("E7").Value = '  -8.47%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  -8.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.22'
$ws.Range("E10").Value = '  -9.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.50'
$ws.Range("E11").Value = '  -3.01%  '
$ws.Range("E12").Value = '  -9.84%  '
$ws.Range("E13").Value = '  -9.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.108'
$ws.Range("E14").Value = '  -1.08%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.870'
$ws.Range("E15").Value = '  -12.05%  '
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '2.582.69'
$ws.Range("E16").Value = '  -5.63%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.48'
$ws.Range("E17").Value = '  -6.71%  '
$ws.Range("D18").Value = '2.238.14'
$ws.Range("E18").Value = '  -5.16%  '
$ws.Range("D19").Value = '43.186.21'
$ws.Range("E19").Value = '  -4.91%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.65'
$ws.Range("E20").Value = '  -3.83%  '
$ws.Range("D21").Value = '0.0₃0972'
$ws.Range("E21").Value = '  -8.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.55'
$ws.Range("E22").Value = '  -10.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.19'
$ws.Range("E24").Value = '  -13.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '238.64'
$ws.Range("E25").Value = '  -8.67%  '
$ws.Range("E26").Value = '  -7.18%  '
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("E28").Value = '  +1.17%  '
$ws.Range("E29").Value = '  -1.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.05'
$ws.Range("E30").Value = '  -10.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.37'
$ws.Range("E31").Value = '  -16.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '36.37'
$ws.Range("E32").Value = '  -3.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0883'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '20.46'
$ws.Range("E34").Value = '  -9.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '150.49'
$ws.Range("E35").Value = '  -10.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.71'
$ws.Range("E36").Value = '  -6.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.21'
$ws.Range("E38").Value = '  +1.41%  '
$ws.Range("E39").Value = '  -7.99%  '
$ws.Range("E40").Value = '  -6.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.105'
$ws.Range("E41").Value = '  -10.47%  '
$ws.Range("E42").Value = '  -8.38%  '
$ws.Range("E43").Value = '  -8.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.68'
$ws.Range("E44").Value = '  +4.79%  '
$ws.Range("E45").Value = '  +0.13%  '
$ws.Range("D46").Value = '1.747.67'
$ws.Range("E46").Value = '  -5.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '86.83'
$ws.Range("E47").Value = '  -10.67%  '
$ws.Range("E48").Value = '  -10.14%  '
$ws.Range("E49").Value = '  -9.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '75.74'
$ws.Range("E50").Value = '  -10.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '59.11'
$ws.Range("E51").Value = '  -16.52%  '
